# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off again (new xliff files produced), which makes its "latest
# handback" file out of date, producing a new status/error for both the
# zh-cn and de-de locales, plus an updated summary on the Overview sheet.

$wb = $excel.ActiveWorkbook

# Helper: assign a text value to a cell while guaranteeing it is stored as
# plain text (t="s"), even for values like "True"/"False" that Excel would
# otherwise coerce into a native boolean type. We do this by writing the
# text through a formula in a scratch cell, converting that formula to a
# plain value via PasteSpecial (values only), and then copying the
# resulting (already-text-typed) value onto the destination cell. Using
# the destination's own .Value setter this way keeps the destination's
# existing number format/style untouched.
function Set-TextValue($Worksheet, $Address, $Text) {
    $scratch = $Worksheet.Range("ZZ1000")
    $escaped = $Text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy($scratch)
    $scratch.PasteSpecial(-4163) # xlPasteValues
    $Worksheet.Range($Address).Value = $scratch.Value
    $scratch.ClearContents()
}

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md file. Both locale status columns move
# to "Ready for handoff" and the latest generate-date timestamp updates.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 12:35:21"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md file.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsZhCn "F3" "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-17 12:35:16"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8481a5a444949f897dce4b93e52bf4db989a638/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e50d556df290eb9f5554e63ff1fafd8a15edcf85/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md file.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsDeDe "F3" "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-17 12:35:21"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8481a5a444949f897dce4b93e52bf4db989a638/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e50d556df290eb9f5554e63ff1fafd8a15edcf85/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
